# Auto-generated edit script: refresh Leve market-price snapshot data
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 54: Family Secrets
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 31254700
$ws.Range("I61").Value = 38465250
$ws.Range("K61").Value = 38465250
$ws.Range("M61").Value = -38465038
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 28613358
$ws.Range("I132").Value = 2087.8667
$ws.Range("K132").Value = 6263.6001
$ws.Range("M132").Value = -3733.6001
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 31254700
$ws.Range("I136").Value = 38465250
$ws.Range("K136").Value = 115395750
$ws.Range("M136").Value = -115393200

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5: Axe Me Anything
$ws.Range("H5").Value = 65.8
$ws.Range("I5").Value = 58
$ws.Range("K5").Value = 58
$ws.Range("M5").Value = 55
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1365159.9
$ws.Range("I86").Value = 2093371.9
$ws.Range("K86").Value = 2093371.9
$ws.Range("M86").Value = -2092248.9
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1365159.9
$ws.Range("I89").Value = 2093371.9
$ws.Range("K89").Value = 10466859.5
$ws.Range("M89").Value = -10461243.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 9715
$ws.Range("I86").Value = 9620
$ws.Range("K86").Value = 9620
$ws.Range("M86").Value = -8497
# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 9715
$ws.Range("I89").Value = 9620
$ws.Range("K89").Value = 48100
$ws.Range("M89").Value = -42484
# Row 94: Beech, Please
$ws.Range("H94").Value = 2060.2307
$ws.Range("I94").Value = 1949.75
$ws.Range("J94").Value = 2080.318
$ws.Range("K94").Value = 1949.75
$ws.Range("L94").Value = 2080.318
$ws.Range("M94").Value = -1498.75
$ws.Range("N94").Value = -2982.318
# Row 99: O Pine
$ws.Range("H99").Value = 21415
$ws.Range("I99").Value = 22499.691
$ws.Range("K99").Value = 22499.691
$ws.Range("M99").Value = -21001.691
# Row 126: A Better Conductor
$ws.Range("H126").Value = 21415
$ws.Range("I126").Value = 22499.691
$ws.Range("K126").Value = 67499.073
$ws.Range("M126").Value = -65029.073

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 109: Cure for What Ails
$ws.Range("H109").Value = 2392.5
$ws.Range("I109").Value = 1997.5
$ws.Range("J109").Value = 2458.3333
$ws.Range("K109").Value = 5992.5
$ws.Range("L109").Value = 7374.999899999999
$ws.Range("M109").Value = -4952.5
$ws.Range("N109").Value = -9454.999899999999
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 2420.2856
$ws.Range("I113").Value = 935.7143
$ws.Range("J113").Value = 3162.5715
$ws.Range("K113").Value = 2807.1429
$ws.Range("L113").Value = 9487.7145
$ws.Range("M113").Value = -637.1428999999998
$ws.Range("N113").Value = -13827.7145
# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 5072.1816
$ws.Range("I134").Value = 2866
$ws.Range("K134").Value = 8598
$ws.Range("M134").Value = -3528
# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 1061
$ws.Range("I140").Value = 702.2105
$ws.Range("K140").Value = 2106.6315
$ws.Range("M140").Value = 3073.3685

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 3454.4092
$ws.Range("I70").Value = 3522.3333
$ws.Range("K70").Value = 3522.3333
$ws.Range("M70").Value = -3252.3333
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 3454.4092
$ws.Range("I73").Value = 3522.3333
$ws.Range("K73").Value = 3522.3333
$ws.Range("M73").Value = -2586.3333
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1545.7059
$ws.Range("I132").Value = 1392.1538
$ws.Range("K132").Value = 4176.4614
$ws.Range("M132").Value = -1646.4614

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 6825
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 6825
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6825
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -7415
# Row 27: Fire and Hide
$ws.Range("H27").Value = 6825
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 6825
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6825
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = -7039
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 5289.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5289.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5289.25
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -5665.25
# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2734.625
$ws.Range("I82").Value = 1520.1111
$ws.Range("J82").Value = 3463.3333
$ws.Range("K82").Value = 1520.1111
$ws.Range("L82").Value = 3463.3333
$ws.Range("M82").Value = -1159.1111
$ws.Range("N82").Value = -4185.3333
# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2734.625
$ws.Range("I85").Value = 1520.1111
$ws.Range("J85").Value = 3463.3333
$ws.Range("K85").Value = 1520.1111
$ws.Range("L85").Value = 3463.3333
$ws.Range("M85").Value = -272.1111000000001
$ws.Range("N85").Value = -5959.3333
# Row 122: Hell on Leather
$ws.Range("H122").Value = 4849.926
$ws.Range("I122").Value = 3846.2
$ws.Range("J122").Value = 5440.353
$ws.Range("K122").Value = 11538.6
$ws.Range("L122").Value = 16321.059
$ws.Range("M122").Value = -9088.599999999999
$ws.Range("N122").Value = -21221.059
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 5657.3
$ws.Range("I136").Value = 3866.6667
$ws.Range("J136").Value = 6424.7144
$ws.Range("K136").Value = 11600.0001
$ws.Range("L136").Value = 19274.1432
$ws.Range("M136").Value = -9050.000100000001
$ws.Range("N136").Value = -24374.1432

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 10: Just for Kecks
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = $null
# Row 27: Hitting Below the Belt
$ws.Range("H27").Value = 73997
$ws.Range("J27").Value = 73997
$ws.Range("L27").Value = 73997
$ws.Range("N27").Value = -74135
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1951.5714
$ws.Range("I81").Value = 1889.75
$ws.Range("J81").Value = 2034
$ws.Range("K81").Value = 3779.5
$ws.Range("L81").Value = 4068
$ws.Range("M81").Value = -2718.5
$ws.Range("N81").Value = -6190
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1951.5714
$ws.Range("I84").Value = 1889.75
$ws.Range("J84").Value = 2034
$ws.Range("K84").Value = 18897.5
$ws.Range("L84").Value = 20340
$ws.Range("M84").Value = -13593.5
$ws.Range("N84").Value = -30948

